# Reorder the monthly rows so that within each year block, October, November
# and December come first, followed by January through September.
#
# The worksheet has a header in row 1 (columns B, C, D) and data rows 2-49
# with column A holding a "YYYY-MM" label and columns B, C, D holding numeric
# values. We read all the existing data keyed by its YYYY-MM label, then
# rewrite the rows in the new order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the extent of the data (header row 1, data starts row 2).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp = -4162
if ($lastRow -lt 2) { $lastRow = 1 }

$firstDataRow = 2

# Read all existing rows into a lookup keyed by the date label in column A.
$data = @{}
$order = @()

for ($r = $firstDataRow; $r -le $lastRow; $r++) {
    $label = $ws.Cells.Item($r, 1).Value()
    if ($null -eq $label -or $label -eq "") { continue }
    $b = $ws.Cells.Item($r, 2).Value()
    $c = $ws.Cells.Item($r, 3).Value()
    $d = $ws.Cells.Item($r, 4).Value()
    $data[$label] = @($b, $c, $d)
    $order += $label
}

# Figure out which years are present, preserving first-seen order.
$years = @()
foreach ($label in $order) {
    $year = $label.Substring(0, 4)
    if (-not ($years -contains $year)) {
        $years += $year
    }
}

# Build the new ordering: for each year, months 10, 11, 12, then 1-9.
$newOrder = @()
foreach ($year in $years) {
    foreach ($m in 10..12) {
        $label = "{0}-{1:D2}" -f $year, $m
        if ($data.ContainsKey($label)) {
            $newOrder += $label
        }
    }
    foreach ($m in 1..9) {
        $label = "{0}-{1:D2}" -f $year, $m
        if ($data.ContainsKey($label)) {
            $newOrder += $label
        }
    }
}

# Write the data back out in the new order.
$r = $firstDataRow
foreach ($label in $newOrder) {
    $vals = $data[$label]
    $ws.Cells.Item($r, 1).Value = $label
    $ws.Cells.Item($r, 2).Value = $vals[0]
    $ws.Cells.Item($r, 3).Value = $vals[1]
    $ws.Cells.Item($r, 4).Value = $vals[2]
    $r++
}
